$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 17 (pushing New Zealand..USA down by one),
# so that "Netherlands" is placed alphabetically between Luxembourg and New Zealand.
$ws.Rows(17).Insert()

# Copy formatting from the row above (Luxembourg, row 16) so the new row matches
# the existing number formatting (percentages) used throughout the table. Only
# copy the used columns (A:G) to avoid bleeding formatting into the rest of the row.
$ws.Range("A16:G16").Copy()
$ws.Range("A17:G17").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = $false

# Tiny floating-point recalculation artifact on the pre-existing Luxembourg row
# (last-bit difference from a source-data refresh upstream of this edit).
$ws.Range("B16").Value = 0.00618825317555097

# Fill in the data for Netherlands.
$ws.Range("A17").Value = "Netherlands"
$ws.Range("B17").Value = 0.08754907717
$ws.Range("C17").Value = 0.10144263489
$ws.Range("D17").Value = 0.13979585933
$ws.Range("E17").Value = 0.08975508656
$ws.Range("F17").Value = 0.0552086109
$ws.Range("G17").Value = 0.04355178694

$wb.Save()
